$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "60.453.76"
$ws.Range("E2").Value = "  +2.05%  "

# Row 3
$ws.Range("D3").Value = "2.610.34"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
Set-TextValue "D5" "563.82"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6
Set-TextValue "D6" "142.51"
$ws.Range("E6").Value = "  -0.70%  "

# Row 7
Set-TextValue "D7" "0.996"
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "2.632.39"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10
$ws.Range("E10").Value = "  +0.38%  "

# Row 11
Set-TextValue "D11" "0.106"
$ws.Range("E11").Value = "  +1.81%  "

# Row 12
$ws.Range("E12").Value = "  +2.65%  "

# Row 13
Set-TextValue "D13" "0.372"
$ws.Range("E13").Value = "  +7.76%  "

# Row 14
$ws.Range("D14").Value = "3.073.93"
$ws.Range("E14").Value = "  +1.03%  "

# Row 15
$ws.Range("D15").Value = "60.384.68"
$ws.Range("E15").Value = "  +1.87%  "

# Row 16
Set-TextValue "D16" "23.42"
$ws.Range("E16").Value = "  +4.00%  "

# Row 17
$ws.Range("E17").Value = "  +1.67%  "

# Row 18
$ws.Range("D18").Value = "2.618.55"
$ws.Range("E18").Value = "  +1.15%  "

# Row 19
Set-TextValue "D19" "4.64"
$ws.Range("E19").Value = "  +2.57%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "10.86"
$ws.Range("E20").Value = "  +6.22%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "344.91"
$ws.Range("E21").Value = "  +2.57%  "

# Row 22
Set-TextValue "D22" "6.90"
$ws.Range("E22").Value = "  +11.46%  "

# Row 24
$ws.Range("E24").Value = "  +13.94%  "

# Row 25
Set-TextValue "D25" "63.12"
$ws.Range("E25").Value = "  -1.55%  "

# Row 26
Set-TextValue "D26" "0.994"
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("E27").Value = "  -0.30%  "

# Row 28
Set-TextValue "D28" "7.72"
$ws.Range("E28").Value = "  +6.42%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0795"
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("E31").Value = "  +2.67%  "

# Row 32
Set-TextValue "D32" "1.71"
$ws.Range("E32").Value = "  +1.37%  "

# Row 33
Set-TextValue "D33" "159.89"
$ws.Range("E33").Value = "  +2.12%  "

# Row 34
$ws.Range("E34").Value = "  +2.33%  "

# Row 35
Set-TextValue "D35" "4.22"
$ws.Range("E35").Value = "  +4.32%  "

# Row 36
Set-TextValue "D36" "0.969"
$ws.Range("E36").Value = "  +8.97%  "

# Row 37
$ws.Range("E37").Value = "  +4.62%  "

# Row 38
Set-TextValue "D38" "1.58"
$ws.Range("E38").Value = "  +5.45%  "

# Row 39
Set-TextValue "D39" "37.73"
$ws.Range("E39").Value = "  +2.47%  "

# Row 40
Set-TextValue "D40" "0.857"
$ws.Range("E40").Value = "  -2.67%  "

# Row 41
Set-TextValue "D41" "3.79"
$ws.Range("E41").Value = "  +3.44%  "

# Row 42
Set-TextValue "D42" "300.13"
$ws.Range("E42").Value = "  +1.59%  "

# Row 43
Set-TextValue "D43" "140.61"
$ws.Range("E43").Value = "  +12.98%  "

# Row 44
Set-TextValue "D44" "0.995"
$ws.Range("E44").Value = "  -0.35%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.0981"
$ws.Range("E45").Value = "  +0.50%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D46" "0.604"
$ws.Range("E46").Value = "  +1.05%  "

# Row 47
Set-TextValue "D47" "0.0243"
$ws.Range("E47").Value = "  +4.49%  "

# Row 48
Set-TextValue "D48" "0.0543"
$ws.Range("E48").Value = "  +0.89%  "

# Row 49
$ws.Range("E49").Value = "  +0.53%  "

# Row 50
Set-TextValue "D50" "19.50"
$ws.Range("E50").Value = "  +5.11%  "

# Row 51
Set-TextValue "D51" "4.80"
$ws.Range("E51").Value = "  +6.59%  "
